$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded. It belongs right after the
# existing row 248 (chronologically the data set is built by prepending the
# newest entry), so insert a fresh row at 249 and shift every following row
# down by one (this also pushes the former last row, 340, down to 341).
$ws.Rows("249:249").Insert()

# Populate the newly inserted row 249 with the new weekly record.
$ws.Cells.Item(249, 1).Value = 10
$ws.Cells.Item(249, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(249, 3).Value = "La Araucanía"
$ws.Cells.Item(249, 4).Value = 45146
$ws.Cells.Item(249, 5).Value = 9
$ws.Cells.Item(249, 6).Value = 100114007
$ws.Cells.Item(249, 7).Value = "Jengibre"
$ws.Cells.Item(249, 8).Value = "Sin especificar"
$ws.Cells.Item(249, 9).Value = "Primera"
$ws.Cells.Item(249, 10).Value = 25
$ws.Cells.Item(249, 11).Value = 22000
$ws.Cells.Item(249, 12).Value = 22000
$ws.Cells.Item(249, 13).Value = 22000
$ws.Cells.Item(249, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(249, 15).Value = "Perú"
$ws.Cells.Item(249, 16).Value = 1692
$ws.Cells.Item(249, 17).Value = 13
$ws.Cells.Item(249, 18).Value = "Hortaliza"
